# "Fruta / hortaliza, semanal" weekly update:
# Insert a new weekly price row at the top of the Brócoli data block
# (row 419), shifting all the existing data rows down by one
# (old row 419 -> 420, ..., old row 440 -> 441) and filling the new
# row 419 with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 419:440 down to 420:441, inserting a blank row at 419
# (mirrors Excel's Rows.Insert, which also copies the row-above's
# formatting onto the new row, e.g. the date number format on column D).
$ws.Rows.Item(419).Insert()

$ws.Range("A419").Value = 1
$ws.Range("B419").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C419").Value = "Arica y Parinacota"
$ws.Range("D419").Value = 44826
$ws.Range("E419").Value = 15
$ws.Range("F419").Value = 100112023
$ws.Range("G419").Value = "Brócoli"
$ws.Range("H419").Value = "Sin especificar"
$ws.Range("I419").Value = "Tercera"
$ws.Range("J419").Value = 1000
$ws.Range("K419").Value = 500
$ws.Range("L419").Value = 600
$ws.Range("M419").Value = 550
$ws.Range("N419").Value = "`$/unidad"
$ws.Range("O419").Value = "Región de Arica y Parinacota"
$ws.Range("P419").Value = 550
$ws.Range("Q419").Value = 1
$ws.Range("R419").Value = "Hortaliza"
